# Apply updated crypto price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store the value as text (matches
# the source data's inline-string / General-format text cells) instead
# of auto-converting numeric-looking strings like '0.510' or '27.618.22'.
$ws.Range("D2").Value = "'" + '27.618.22'
$ws.Range("E2").Value = "'" + '  -2.43%  '
$ws.Range("D3").Value = "'" + '1.659.22'
$ws.Range("E3").Value = "'" + '  -4.21%  '
$ws.Range("E4").Value = "'" + '  -0.07%  '
$ws.Range("D5").Value = "'" + '215.14'
$ws.Range("E5").Value = "'" + '  -1.95%  '
$ws.Range("D6").Value = "'" + '0.510'
$ws.Range("E6").Value = "'" + '  -2.29%  '
$ws.Range("E7").Value = "'" + '  -0.04%  '
$ws.Range("D8").Value = "'" + '24.18'
$ws.Range("E8").Value = "'" + '  +0.19%  '
$ws.Range("D9").Value = "'" + '0.262'
$ws.Range("E9").Value = "'" + '  -2.35%  '
$ws.Range("D10").Value = "'" + '0.0620'
$ws.Range("E10").Value = "'" + '  -2.67%  '
$ws.Range("D11").Value = "'" + '0.0880'
$ws.Range("E11").Value = "'" + '  -1.68%  '
$ws.Range("D12").Value = "'" + '1.894.09'
$ws.Range("E12").Value = "'" + '  -4.21%  '
$ws.Range("D13").Value = "'" + '1.644.37'
$ws.Range("E13").Value = "'" + '  -5.11%  '
$ws.Range("E14").Value = "'" + '  -2.83%  '
$ws.Range("E15").Value = "'" + '  +0.21%  '
$ws.Range("D16").Value = "'" + '65.94'
$ws.Range("E16").Value = "'" + '  -2.76%  '
$ws.Range("D17").Value = "'" + '27.592.64'
$ws.Range("E17").Value = "'" + '  -2.53%  '
$ws.Range("D18").Value = "'" + '241.59'
$ws.Range("E18").Value = "'" + '  -1.08%  '
$ws.Range("E19").Value = "'" + '  -3.42%  '
$ws.Range("D20").Value = "'" + '7.58'
$ws.Range("E20").Value = "'" + '  -4.69%  '
$ws.Range("E21").Value = "'" + '  -0.03%  '
$ws.Range("D22").Value = "'" + '4.48'
$ws.Range("E22").Value = "'" + '  -3.77%  '
$ws.Range("E23").Value = "'" + '  -3.47%  '
$ws.Range("E24").Value = "'" + '  -2.38%  '
$ws.Range("D25").Value = "'" + '146.07'
$ws.Range("E25").Value = "'" + '  -2.22%  '
$ws.Range("D26").Value = "'" + '7.21'
$ws.Range("E26").Value = "'" + '  -4.31%  '
$ws.Range("D27").Value = "'" + '16.32'
$ws.Range("E27").Value = "'" + '  -2.09%  '
$ws.Range("E28").Value = "'" + '  -0.12%  '
$ws.Range("D29").Value = "'" + '0.112'
$ws.Range("E29").Value = "'" + '  -2.33%  '
$ws.Range("E30").Value = "'" + '  +0.58%  '
$ws.Range("E31").Value = "'" + '  -2.89%  '
$ws.Range("E32").Value = "'" + '  -2.74%  '
$ws.Range("D33").Value = "'" + '1.456.39'
$ws.Range("E33").Value = "'" + '  -2.30%  '
$ws.Range("E34").Value = "'" + '  -5.09%  '
$ws.Range("E35").Value = "'" + '  -4.94%  '
$ws.Range("E36").Value = "'" + '  -1.22%  '
$ws.Range("D37").Value = "'" + '0.926'
$ws.Range("E37").Value = "'" + '  -5.18%  '
$ws.Range("E38").Value = "'" + '  -2.48%  '
$ws.Range("E39").Value = "'" + '  -5.03%  '
$ws.Range("D40").Value = "'" + '69.81'
$ws.Range("E40").Value = "'" + '  -0.63%  '
$ws.Range("E41").Value = "'" + '  -4.29%  '
$ws.Range("E42").Value = "'" + '  -0.06%  '
$ws.Range("E43").Value = "'" + '  -4.39%  '
$ws.Range("E44").Value = "'" + '  -3.56%  '
$ws.Range("D45").Value = "'" + '0.793'
$ws.Range("E45").Value = "'" + '  -0.33%  '
$ws.Range("D46").Value = "'" + '1.800.97'
$ws.Range("E46").Value = "'" + '  -4.15%  '
$ws.Range("E47").Value = "'" + '  -0.94%  '
$ws.Range("D48").Value = "'" + '88.73'
$ws.Range("E48").Value = "'" + '  -2.31%  '
$ws.Range("D49").Value = "'" + '0.0₆0107'
$ws.Range("E49").Value = "'" + '  -5.77%  '
$ws.Range("E50").Value = "'" + '  -1.64%  '
$ws.Range("D51").Value = "'" + '7.84'
$ws.Range("E51").Value = "'" + '  -4.27%  '
